$wb = $excel.ActiveWorkbook

# ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 107.29
$ws.Range("I15").Value = 107.29
$ws.Range("K15").Value = 321.87
$ws.Range("M15").Value = -152.87

# ALC row 26
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("N26").ClearContents()

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 764.3913
$ws.Range("I33").Value = 598.95
$ws.Range("J33").Value = 1867.3334
$ws.Range("K33").Value = 598.95
$ws.Range("L33").Value = 1867.3334
$ws.Range("M33").Value = -369.95
$ws.Range("N33").Value = -2325.3334

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10510.333
$ws.Range("I62").Value = 12598.833
$ws.Range("J62").Value = 6333.3335
$ws.Range("K62").Value = 12598.833
$ws.Range("L62").Value = 6333.3335
$ws.Range("M62").Value = -11974.833
$ws.Range("N62").Value = -7581.3335

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10510.333
$ws.Range("I65").Value = 12598.833
$ws.Range("J65").Value = 6333.3335
$ws.Range("K65").Value = 62994.165
$ws.Range("L65").Value = 31666.6675
$ws.Range("M65").Value = -59874.165
$ws.Range("N65").Value = -37906.6675

# ALC row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 529.5769
$ws.Range("I101").Value = 412.66666
$ws.Range("J101").Value = 591.4706
$ws.Range("K101").Value = 1237.99998
$ws.Range("L101").Value = 1774.4118
$ws.Range("M101").Value = 384.0000199999999
$ws.Range("N101").Value = -5018.4118

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 53829.816
$ws.Range("I116").Value = 75454.39999999999
$ws.Range("J116").Value = 7491.4287
$ws.Range("K116").Value = 75454.39999999999
$ws.Range("L116").Value = 7491.4287
$ws.Range("M116").Value = -72012.39999999999
$ws.Range("N116").Value = -14375.4287

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4883.114
$ws.Range("I132").Value = 3214.5483
$ws.Range("K132").Value = 9643.644899999999
$ws.Range("M132").Value = -7113.644899999999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4361.2095
$ws.Range("I138").Value = 2706.8572
$ws.Range("J138").Value = 4843.729
$ws.Range("K138").Value = 8120.571599999999
$ws.Range("L138").Value = 14531.187
$ws.Range("M138").Value = -2980.571599999999
$ws.Range("N138").Value = -24811.187

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26115.918
$ws.Range("I32").Value = 14303.708
$ws.Range("J32").Value = 47923.08
$ws.Range("K32").Value = 14303.708
$ws.Range("L32").Value = 47923.08
$ws.Range("M32").Value = -14016.708
$ws.Range("N32").Value = -48497.08

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 48652.305
$ws.Range("I134").Value = 102180.3
$ws.Range("J134").Value = 7476.923
$ws.Range("K134").Value = 306540.9
$ws.Range("L134").Value = 22430.769
$ws.Range("M134").Value = -304005.9
$ws.Range("N134").Value = -27500.769

# CRP row 17
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 3000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# CRP row 25
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 4000
$ws.Range("K25").Value = 4000
$ws.Range("M25").Value = -3826

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2994.0364
$ws.Range("I31").Value = 2238.7368
$ws.Range("J31").Value = 4682.353
$ws.Range("K31").Value = 2238.7368
$ws.Range("L31").Value = 4682.353
$ws.Range("M31").Value = -1943.7368
$ws.Range("N31").Value = -5272.353

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2994.0364
$ws.Range("I34").Value = 2238.7368
$ws.Range("J34").Value = 4682.353
$ws.Range("K34").Value = 2238.7368
$ws.Range("L34").Value = 4682.353
$ws.Range("M34").Value = -2036.7368
$ws.Range("N34").Value = -5086.353

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4763.8716
$ws.Range("I58").Value = 5072.1035
$ws.Range("J58").Value = 3870
$ws.Range("K58").Value = 5072.1035
$ws.Range("L58").Value = 3870
$ws.Range("M58").Value = -4869.1035
$ws.Range("N58").Value = -4276

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 86217.664
$ws.Range("I99").Value = 145687.42
$ws.Range("J99").Value = 2960
$ws.Range("K99").Value = 145687.42
$ws.Range("L99").Value = 2960
$ws.Range("M99").Value = -144189.42
$ws.Range("N99").Value = -5956

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 86217.664
$ws.Range("I126").Value = 145687.42
$ws.Range("J126").Value = 2960
$ws.Range("K126").Value = 437062.26
$ws.Range("L126").Value = 8880
$ws.Range("M126").Value = -434592.26
$ws.Range("N126").Value = -13820

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4763.8716
$ws.Range("I136").Value = 5072.1035
$ws.Range("J136").Value = 3870
$ws.Range("K136").Value = 15216.3105
$ws.Range("L136").Value = 11610
$ws.Range("M136").Value = -12666.3105
$ws.Range("N136").Value = -16710

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 705.54285
$ws.Range("J5").Value = 1073
$ws.Range("L5").Value = 3219
$ws.Range("N5").Value = -3443

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4255.636
$ws.Range("I70").Value = 1228
$ws.Range("J70").Value = 5985.7144
$ws.Range("K70").Value = 3684
$ws.Range("L70").Value = 17957.1432
$ws.Range("M70").Value = -3369
$ws.Range("N70").Value = -18587.1432

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 4255.636
$ws.Range("I73").Value = 1228
$ws.Range("J73").Value = 5985.7144
$ws.Range("K73").Value = 3684
$ws.Range("L73").Value = 17957.1432
$ws.Range("M73").Value = -2592
$ws.Range("N73").Value = -20141.1432

# CUL row 87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11335.667
$ws.Range("I87").Value = 12007
$ws.Range("K87").Value = 36021
$ws.Range("M87").Value = -34773

# CUL row 90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 11335.667
$ws.Range("I90").Value = 12007
$ws.Range("K90").Value = 108063
$ws.Range("M90").Value = -101823

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 705.54285
$ws.Range("J135").Value = 1073
$ws.Range("L135").Value = 9657
$ws.Range("N135").Value = -14727

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5408.1665
$ws.Range("I102").Value = 6334.385
$ws.Range("K102").Value = 6334.385
$ws.Range("M102").Value = -4712.385

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 325.17392
$ws.Range("I55").Value = 358.30768
$ws.Range("J55").Value = 282.1
$ws.Range("K55").Value = 358.30768
$ws.Range("L55").Value = 282.1
$ws.Range("M55").Value = -185.30768
$ws.Range("N55").Value = -628.1

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5577.5835
$ws.Range("I136").Value = 3237.0952
$ws.Range("J136").Value = 8854.267
$ws.Range("K136").Value = 9711.285600000001
$ws.Range("L136").Value = 26562.801
$ws.Range("M136").Value = -7161.285600000001
$ws.Range("N136").Value = -31662.801

# WVR row 6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("N6").ClearContents()

# WVR row 46
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 59028.5
$ws.Range("J46").Value = 59028.5
$ws.Range("L46").Value = 59028.5
$ws.Range("N46").Value = -59490.5

# WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3169.9285
$ws.Range("I96").Value = 1485.8572
$ws.Range("J96").Value = 4854
$ws.Range("K96").Value = 1485.8572
$ws.Range("L96").Value = 4854
$ws.Range("M96").Value = -112.8571999999999
$ws.Range("N96").Value = -7600

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2780.743
$ws.Range("I132").Value = 1705.8182
$ws.Range("K132").Value = 5117.4546
$ws.Range("M132").Value = -2587.4546

# WVR row 134
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 59028.5
$ws.Range("J134").Value = 59028.5
$ws.Range("L134").Value = 177085.5
$ws.Range("N134").Value = -182155.5
Write-Host "Applied scheduled runner price/profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets."
